# PlayerPerformance_5661.xlsx edit
# 1. Clear the (spurious, empty) INNING_NUMBER cells B2, B3, B4, B6 on the
#    "ODI Batting" sheet so they become true blanks again.
# 2. Add a new "ODI Batting Extra" sheet (after "ODI Bowling") with
#    MATCH_CODE / BATTING_POSITION / NUM_4 / NUM_6 / PERCENT_RUNS_OF_TOTAL /
#    MAN_OF_MATCH columns, matching the header style used on the other
#    sheets, and populate it with the per-match extra-batting data.

$wb = $excel.ActiveWorkbook

# --- 1. Clean up stray empty cells on "ODI Batting" ------------------------
$batting = $wb.Worksheets.Item("ODI Batting")
$batting.Range("B2").ClearContents()
$batting.Range("B3").ClearContents()
$batting.Range("B4").ClearContents()
$batting.Range("B6").ClearContents()

# E3 held a stray non-breaking space; normalize to a regular space.
$batting.Range("E3").Value = " "

# --- 2. Create the new "ODI Batting Extra" sheet ----------------------------
$bowling = $wb.Worksheets.Item("ODI Bowling")
$extra = $wb.Worksheets.Add($null, $bowling)
$extra.Name = "ODI Batting Extra"

# Header row (copy the bold/border/centered formatting used on every other
# sheet's header row, then overwrite the text).
$playerInfo = $wb.Worksheets.Item("Player Info")
$playerInfo.Range("A1:D1").Copy()
$extra.Range("A1:F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$extra.Range("A1").Value = "MATCH_CODE"
$extra.Range("B1").Value = "BATTING_POSITION"
$extra.Range("C1").Value = "NUM_4"
$extra.Range("D1").Value = "NUM_6"
$extra.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$extra.Range("F1").Value = "MAN_OF_MATCH"

# MATCH_CODE (col A) and the textual NUM_4/NUM_6/PERCENT columns on row 4
# hold digit-only / percent-looking strings that must stay text, not get
# auto-coerced to numbers -- force a text format before writing them.
$extra.Range("A2:A6").NumberFormat = "@"
$extra.Range("C4:E4").NumberFormat = "@"

# Data rows
$extra.Range("A2").Value = "4401"
$extra.Range("F2").Value = "NO"

$extra.Range("A3").Value = "4405"
$extra.Range("F3").Value = "NO"

$extra.Range("A4").Value = "4472"
$extra.Range("B4").Value = 11
$extra.Range("C4").Value = "1"
$extra.Range("D4").Value = "0"
$extra.Range("E4").Value = "2.83%"
$extra.Range("F4").Value = "NO"

$extra.Range("A5").Value = "4473"

$extra.Range("A6").Value = "4476"
